$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data rows 11-14 (17-class model results + ensemble code) ---
$newRows = @(
    @(220623, "문현우", "EffNetV2M",          17, 5, 0.8951, 0.8888,               0.4343, 0.5318),
    @(220623, "박영서", "Xception",           17, 5, 0.8742, 0.874048471450805,    0.7677, 0.75166046619415205),
    @(220623, "박영서", "VGG16",              17, 5, 0.78,   0.76257210969924905,  0.9666, 1.3760806322097701),
    @(220623, "임한준", "InceptionResNetV2",  17, 5, 0.8371, 0.824221432209014,    0.7708, 0.89661860466003396)
)

$r = 11
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

# H12 keeps the header/general style (same as H1) rather than the percent-format
# style used by the other H-column cells - matches the source edit exactly.
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(12, 8).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- AutoFilter over the full data range ---
$ws.Range("A1:I14").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$I`$14")
$filterName.Visible = $false

# --- Conditional formatting (3-color scale) on F2:G1048576 ---
$cfRange = $ws.Range("F2:G1048576")
$cfRange.FormatConditions.AddColorScale(3)

# --- View: freeze header row, zoom to 115%, selection on C12 ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 115
$ws.Range("C12").Select()
